$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" column header, styled like the other header cells (column G1 "sum")
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# New "Save" column values for rows 2-10
$saveValues = @(0, 1, 1, 1, 0, 0, 1, 1, 1)
for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
